$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.5
$ws.Range("I2").Value = 2.9
$ws.Range("J2").Value = 3.2
$ws.Range("N2").Value = 9
$ws.Range("X2").Value = 12
$ws.Range("AC2").Value = 9
$ws.Range("AH2").Value = 13
$ws.Range("AK2").Value = 23
$ws.Range("AM2").Value = 4.5
$ws.Range("AN2").Value = 15
$ws.Range("AP2").Value = 51
